# Adds a 9th "leave" column (I) to the roster sheet:
#  - I1 header mirrors F1/G1/H1 ("leave")
#  - I2 gets a new leave record for the row-2 employee (Akram Khan)
#  - column I width + row 2 height are adjusted to fit the new content
#  - selection moves to I2, matching the saved state in the source file

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell I1: same text + formatting as the existing "leave" headers ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "leave"

# --- Data cell I2: new leave note, formatted like the other leave-note cells ---
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)   # xlPasteFormats
$newLine = [char]10
$ws.Range("I2").Value = "leaveType:C;" + $newLine + "start: 12 Feb 2020;" + $newLine + "end: 24 Feb 2020;" + $newLine + "specialDays: -1;"

# --- Column width for the new column I ---
$ws.Columns.Item(9).ColumnWidth = 20.33

# --- Row 2 needs to grow to accommodate the extra wrapped text ---
$ws.Rows.Item(2).RowHeight = 95

# --- Clear the clipboard marquee and move the active selection to I2 ---
$excel.CutCopyMode = 0
$ws.Range("I2").Select()
